# Regenerate merged AHB files
# - Rename the diff-table headers from "_old" / "_new" suffixes to the
#   explicit format-version suffixes "_FV2310" / "_FV2404".
# - Turn the data range into a native Excel Table (Table1).
# - Freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into an Excel Table (ListObject), headers already
# present on row 1, so tell Excel to use them (xlYes = 1).
$dataRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $dataRange, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
